$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.2980729807298073
$ws.Range("I3").Value = 0.5506984387838948
$ws.Range("K3").Value = 2686.3

$ws.Range("Q3").Value = 2649
$ws.Range("R3").Value = 2657
$ws.Range("S3").Value = 2665
$ws.Range("T3").Value = 2676
$ws.Range("U3").Value = 2684

$ws.Range("V3").Value = 2219
$ws.Range("W3").Value = 2211
$ws.Range("X3").Value = 2203
$ws.Range("Y3").Value = 2192
$ws.Range("Z3").Value = 2184

$ws.Range("AF3").Value = 0.455834
$ws.Range("AG3").Value = 0.454191
$ws.Range("AH3").Value = 0.452547
$ws.Range("AI3").Value = 0.450288
$ws.Range("AJ3").Value = 0.448644
